# PZG_CelPracy.xlsx - align value texts with the SD XML schema, drop the
# obsolete "code name" rows (23-43) and rename the sheet/tab.
# Commit message: "Walidacja pikow XML zgodnie ze schematem SD"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the wording of several "value" entries (column B). The
#     calculated column C (Kolumna1) re-evaluates automatically since it
#     is a table formula referencing [value] for "this row".
$ws.Range("B5").Value  = "utworzenie bazy danych geodezyjnej ewidencji sieci uzbrojenia terenu (GESUT)"
$ws.Range("B6").Value  = "aktualizacja bazy danych geodezyjnej ewidencji sieci uzbrojenia terenu (GESUT)"
$ws.Range("B7").Value  = "utworzenie bazy danych obiektów topograficznych o szczegółowości zapewniającej tworzenie standardowych opracowań kartograficznych w skalach 1:500-1:5000 (BDOT500)"
$ws.Range("B8").Value  = "aktualizacja bazy danych obiektów topograficznych o szczegółowości zapewniającej tworzenie standardowych opracowań kartograficznych w skalach 1:500-1:5000 (BDOT500)"
$ws.Range("B18").Value = "wznowienie znaków granicznych/wyznaczenie punktów/ustalenie przebiegu granic działek ewidencyjnych"
$ws.Range("B20").Value = "geodezyjna inwentaryzacja obiektów budowlanych"
$ws.Range("B21").Value = "wytyczenie obiektów budowlanych"

# --- Drop the trailing "code name" rows (23-43); the table range,
#     dimension and autofilter all shrink to A1:C22 automatically.
$ws.Range("A23:C43").EntireRow.Delete()

# --- Reset the view: scroll back to the top-left and select the
#     calculated "Kolumna1" column for the now-smaller table.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C2:C22").Select() | Out-Null

# --- Rename the worksheet/tab.
$ws.Name = "Arkusz1"
